# Insert a new data row at row 265 (shifting the existing rows 265-296 down
# to 266-297) and populate the newly inserted row with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(265).Insert()

$ws.Cells.Item(265, 1).Value = 8
$ws.Cells.Item(265, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(265, 3).Value = "Coquimbo"
$ws.Cells.Item(265, 4).Value = 44984
$ws.Cells.Item(265, 5).Value = 4
$ws.Cells.Item(265, 6).Value = 100112037
$ws.Cells.Item(265, 7).Value = "Cebollín"
$ws.Cells.Item(265, 8).Value = "Sin especificar"
$ws.Cells.Item(265, 9).Value = "Primera"
$ws.Cells.Item(265, 10).Value = 1200
$ws.Cells.Item(265, 11).Value = 1200
$ws.Cells.Item(265, 12).Value = 1400
$ws.Cells.Item(265, 13).Value = 1300
$ws.Cells.Item(265, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(265, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(265, 16).Value = 217
$ws.Cells.Item(265, 17).Value = 6
$ws.Cells.Item(265, 18).Value = "Hortaliza"
